$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet tab (was "Through 2022-11-13" -> "Through 2022-11-14")
$ws.Name = "Through 2022-11-14"

# Update the column header text in I1 (shared string) to reflect new date
$ws.Range("I1").Value = "2022 (through 11-14)"

# Update November (row 12) 2022 value
$ws.Range("I12").Value = 47

# Update Total (row 14) 2022 value
$ws.Range("I14").Value = 1445
